$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 973.625
$ws.Range("J17").Value = 1005.2
$ws.Range("L17").Value = 3015.6
$ws.Range("N17").Value = -3351.6

$ws.Range("H33").Value = 487.33334
$ws.Range("I33").Value = 505.83334
$ws.Range("J33").Value = 413.33334
$ws.Range("K33").Value = 505.83334
$ws.Range("L33").Value = 413.33334
$ws.Range("M33").Value = -276.83334
$ws.Range("N33").Value = -871.33334

$ws.Range("H74").Value = 5332.6665
$ws.Range("I74").Value = 5332.6665
$ws.Range("K74").Value = 5332.6665
$ws.Range("M74").Value = -4396.6665

$ws.Range("H77").Value = 5332.6665
$ws.Range("I77").Value = 5332.6665
$ws.Range("K77").Value = 26663.3325
$ws.Range("M77").Value = -21983.3325

$ws.Range("H96").Value = 2097.889
$ws.Range("I96").Value = 2407.4285
$ws.Range("J96").Value = 1014.5
$ws.Range("K96").Value = 7222.2855
$ws.Range("L96").Value = 3043.5
$ws.Range("M96").Value = -5849.2855
$ws.Range("N96").Value = -5789.5

$ws.Range("H111").Value = 3007.25
$ws.Range("I111").Value = 3007.25
$ws.Range("K111").Value = 9021.75
$ws.Range("M111").Value = -5954.75

$ws.Range("H125").Value = 2554.3076
$ws.Range("J125").Value = 2490
$ws.Range("L125").Value = 22410
$ws.Range("N125").Value = -27330

$ws.Range("H132").Value = 10109944
$ws.Range("I132").Value = 12826574
$ws.Range("K132").Value = 38479722
$ws.Range("M132").Value = -38477192

$ws.Range("H137").Value = 3737.9333
$ws.Range("I137").Value = 3996.5
$ws.Range("J137").Value = 3442.4285
$ws.Range("K137").Value = 11989.5
$ws.Range("L137").Value = 10327.2855
$ws.Range("M137").Value = -9439.5
$ws.Range("N137").Value = -15427.2855

$ws.Range("H138").Value = 2731.8489
$ws.Range("I138").Value = 1840.6666
$ws.Range("J138").Value = 2876.3647
$ws.Range("K138").Value = 5521.9998
$ws.Range("L138").Value = 8629.0941
$ws.Range("M138").Value = -381.9997999999996
$ws.Range("N138").Value = -18909.0941

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10776.59
$ws.Range("I32").Value = 7919.9395
$ws.Range("J32").Value = 16321.853
$ws.Range("K32").Value = 7919.9395
$ws.Range("L32").Value = 16321.853
$ws.Range("M32").Value = -7632.9395
$ws.Range("N32").Value = -16895.853

$ws.Range("H74").Value = 1460.1818
$ws.Range("I74").Value = 830.7619
$ws.Range("J74").Value = 2561.6667
$ws.Range("K74").Value = 830.7619
$ws.Range("L74").Value = 2561.6667
$ws.Range("M74").Value = 43.23810000000003
$ws.Range("N74").Value = -4309.6667

$ws.Range("H77").Value = 1460.1818
$ws.Range("I77").Value = 830.7619
$ws.Range("J77").Value = 2561.6667
$ws.Range("K77").Value = 4153.809499999999
$ws.Range("L77").Value = 12808.3335
$ws.Range("M77").Value = 214.1905000000006
$ws.Range("N77").Value = -21544.3335

$ws.Range("H102").Value = 16669132
$ws.Range("I102").Value = 20833914
$ws.Range("J102").Value = 10000
$ws.Range("K102").Value = 20833914
$ws.Range("L102").Value = 10000
$ws.Range("M102").Value = -20832292
$ws.Range("N102").Value = -13244

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3298.7
$ws.Range("I20").Value = 3129.7144
$ws.Range("J20").Value = 3693
$ws.Range("K20").Value = 3129.7144
$ws.Range("L20").Value = 3693
$ws.Range("M20").Value = -2882.7144
$ws.Range("N20").Value = -4187

$ws.Range("H86").Value = 47623176
$ws.Range("I86").Value = 52635668
$ws.Range("K86").Value = 52635668
$ws.Range("M86").Value = -52634545

$ws.Range("H89").Value = 47623176
$ws.Range("I89").Value = 52635668
$ws.Range("K89").Value = 263178340
$ws.Range("M89").Value = -263172724

$ws.Range("H94").Value = 7813180
$ws.Range("I94").Value = 9616024
$ws.Range("K94").Value = 9616024
$ws.Range("M94").Value = -9615573

$ws.Range("H105").Value = 72135740
$ws.Range("I105").Value = 77684500
$ws.Range("J105").Value = 2000
$ws.Range("K105").Value = 77684500
$ws.Range("L105").Value = 2000
$ws.Range("M105").Value = -77682753
$ws.Range("N105").Value = -5494

$ws.Range("H132").Value = 38500
$ws.Range("J132").Value = 38500
$ws.Range("L132").Value = 38500
$ws.Range("N132").Value = -48620

$ws.Range("H134").Value = 4989.423
$ws.Range("I134").Value = 618.3684
$ws.Range("J134").Value = 16853.715
$ws.Range("K134").Value = 1855.1052
$ws.Range("L134").Value = 50561.145
$ws.Range("M134").Value = 679.8948
$ws.Range("N134").Value = -55631.145

$ws.Range("H140").Value = 39188.777
$ws.Range("J140").Value = 39188.777
$ws.Range("L140").Value = 39188.777
$ws.Range("N140").Value = -49548.777

$ws.Range("H141").Value = 99500
$ws.Range("J141").Value = 99500
$ws.Range("L141").Value = 99500
$ws.Range("N141").Value = -109860

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1948.0233
$ws.Range("I31").Value = 1875.5952
$ws.Range("K31").Value = 1875.5952
$ws.Range("M31").Value = -1580.5952

$ws.Range("H34").Value = 1948.0233
$ws.Range("I34").Value = 1875.5952
$ws.Range("K34").Value = 1875.5952
$ws.Range("M34").Value = -1673.5952

$ws.Range("H99").Value = 1463485
$ws.Range("I99").Value = 4387288.5
$ws.Range("J99").Value = 1583.25
$ws.Range("K99").Value = 4387288.5
$ws.Range("L99").Value = 1583.25
$ws.Range("M99").Value = -4385790.5
$ws.Range("N99").Value = -4579.25

$ws.Range("H126").Value = 1463485
$ws.Range("I126").Value = 4387288.5
$ws.Range("J126").Value = 1583.25
$ws.Range("K126").Value = 13161865.5
$ws.Range("L126").Value = 4749.75
$ws.Range("M126").Value = -13159395.5
$ws.Range("N126").Value = -9689.75

$ws.Range("H132").Value = 1987.2972
$ws.Range("I132").Value = 1693.7931
$ws.Range("J132").Value = 3051.25
$ws.Range("K132").Value = 5081.379300000001
$ws.Range("L132").Value = 9153.75
$ws.Range("M132").Value = -2551.379300000001
$ws.Range("N132").Value = -14213.75

$ws.Range("H134").Value = 10639697
$ws.Range("I134").Value = 1429.6129
$ws.Range("J134").Value = 31251342
$ws.Range("K134").Value = 4288.8387
$ws.Range("L134").Value = 93754026
$ws.Range("M134").Value = -1753.8387
$ws.Range("N134").Value = -93759096

$ws.Range("H141").Value = 423579.66
$ws.Range("I141").Value = 13499
$ws.Range("J141").Value = 452871.16
$ws.Range("K141").Value = 13499
$ws.Range("L141").Value = 452871.16
$ws.Range("N141").Value = -463231.16
$ws.Range("M141").Value = -8319

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1606.1818
$ws.Range("J5").Value = 1908.75
$ws.Range("L5").Value = 5726.25
$ws.Range("N5").Value = -5950.25

$ws.Range("H17").Value = 216
$ws.Range("J17").Value = 190
$ws.Range("L17").Value = 570
$ws.Range("N17").Value = -908

$ws.Range("H131").Value = 27820448
$ws.Range("J131").Value = 52851.965
$ws.Range("L131").Value = 158555.895
$ws.Range("N131").Value = -168635.895

$ws.Range("H135").Value = 1606.1818
$ws.Range("J135").Value = 1908.75
$ws.Range("L135").Value = 17178.75
$ws.Range("N135").Value = -22248.75

$ws.Range("H137").Value = 30007168
$ws.Range("I137").Value = 75002670
$ws.Range("J137").Value = 10164.333
$ws.Range("K137").Value = 225008010
$ws.Range("L137").Value = 30492.999
$ws.Range("M137").Value = -225002910
$ws.Range("N137").Value = -40692.999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 471.69565
$ws.Range("I2").Value = 382.15384
$ws.Range("J2").Value = 588.1
$ws.Range("K2").Value = 382.15384
$ws.Range("L2").Value = 588.1
$ws.Range("M2").Value = -269.15384
$ws.Range("N2").Value = -814.1

$ws.Range("H36").Value = 3005.6667
$ws.Range("I36").Value = 3508.5
$ws.Range("J36").Value = 2000
$ws.Range("K36").Value = 3508.5
$ws.Range("L36").Value = 2000
$ws.Range("M36").Value = -3023.5
$ws.Range("N36").Value = -2970

$ws.Range("H134").Value = 27532.8
$ws.Range("J134").Value = 27532.8
$ws.Range("L134").Value = 82598.39999999999
$ws.Range("N134").Value = -87668.39999999999

$ws.Range("H135").Value = 40536.875
$ws.Range("J135").Value = 39185
$ws.Range("L135").Value = 39185
$ws.Range("N135").Value = -49325

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1600.1818
$ws.Range("I16").Value = 1600.1818
$ws.Range("K16").Value = 1600.1818
$ws.Range("M16").Value = -1430.1818

$ws.Range("H22").Value = 1980.6
$ws.Range("J22").Value = 1225.5
$ws.Range("L22").Value = 1225.5
$ws.Range("N22").Value = -1815.5

$ws.Range("H27").Value = 1980.6
$ws.Range("J27").Value = 1225.5
$ws.Range("L27").Value = 1225.5
$ws.Range("N27").Value = -1439.5

$ws.Range("H40").Value = 2296.2778
$ws.Range("I40").Value = 2194.923
$ws.Range("J40").Value = 2559.8
$ws.Range("K40").Value = 2194.923
$ws.Range("L40").Value = 2559.8
$ws.Range("M40").Value = -2058.923
$ws.Range("N40").Value = -2831.8

$ws.Range("H132").Value = 2129.2727
$ws.Range("I132").Value = 1627.625
$ws.Range("J132").Value = 3467
$ws.Range("K132").Value = 4882.875
$ws.Range("L132").Value = 10401
$ws.Range("M132").Value = -2352.875
$ws.Range("N132").Value = -15461

$ws.Range("H136").Value = 1462.4375
$ws.Range("I136").Value = 1359.9333
$ws.Range("K136").Value = 4079.7999
$ws.Range("M136").Value = -1529.7999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 294.92307
$ws.Range("I100").Value = 309.1
$ws.Range("K100").Value = 618.2
$ws.Range("M100").Value = -77.20000000000005

$ws.Range("H107").Value = 349.45456
$ws.Range("J107").Value = 482.125
$ws.Range("L107").Value = 1446.375
$ws.Range("N107").Value = -5286.375

$ws.Range("H132").Value = 3228.8958
$ws.Range("I132").Value = 3274.75
$ws.Range("K132").Value = 9824.25
$ws.Range("M132").Value = -7294.25
